$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 68, shifting existing rows 68-93 down to 69-94.
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row 68 with the new weekly price record.
$ws.Cells.Item(68, 1).Value = 1
$ws.Cells.Item(68, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(68, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(68, 4).Value = 44627
$ws.Cells.Item(68, 5).Value = 15
$ws.Cells.Item(68, 6).Value = 100112008
$ws.Cells.Item(68, 7).Value = "Coliflor"
$ws.Cells.Item(68, 8).Value = "Sin especificar"
$ws.Cells.Item(68, 9).Value = "Segunda"
$ws.Cells.Item(68, 10).Value = 800
$ws.Cells.Item(68, 11).Value = 900
$ws.Cells.Item(68, 12).Value = 1000
$ws.Cells.Item(68, 13).Value = 950
$ws.Cells.Item(68, 14).Value = "$/unidad"
$ws.Cells.Item(68, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(68, 16).Value = 950
$ws.Cells.Item(68, 17).Value = 1
$ws.Cells.Item(68, 18).Value = "Hortaliza"
